$d = $word.ActiveDocument

# 1. Insert a new sub-bullet paragraph "The appropriateness of guerilla
#    activity within an organization." right after the "Ethical obligations
#    as decision-making criteria." bullet, matching its style/level.
$idx = 0
$anchorIdx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Ethical obligations as decision-making criteria*") {
        $anchorIdx = $idx
        break
    }
}

$anchorPara = $d.Paragraphs($anchorIdx)
$anchorPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($anchorIdx + 1)
$newPara.Range.Text = "The appropriateness of guerilla activity within an organization."

# 2. Move the "_GoBack" bookmark from the end of "How authority derives
#    legitimacy." to the start of the "Concluding remarks" paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Concluding remarks*") {
        $r = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $r)
        break
    }
}
